# Regenerate save_data to use K instead of Strike#: recalc and write s_vals
# into column G (header "K") for rows 2-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 4
    6  = 3
    7  = 1
    8  = 3
    9  = 2
    10 = 1
    11 = 0
    12 = 4
    13 = 0
    14 = 4
    15 = 1
    16 = 1
    17 = 0
    18 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
